$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.950.11'
$ws.Range('E2').Value = '  +3.25%  '
$ws.Range('D3').Value = '1.725.02'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.61'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.524'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.95'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +11.91%  '
$ws.Range('E9').Value = '  +3.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0633'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.84%  '
$ws.Range('E11').Value = '  +1.93%  '
$ws.Range('D12').Value = '1.969.54'
$ws.Range('E12').Value = '  +2.99%  '
$ws.Range('D13').Value = '1.728.64'
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('E15').Value = '  +5.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.92'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.52%  '
$ws.Range('D17').Value = '27.901.48'
$ws.Range('E17').Value = '  +3.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.36'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.39%  '
$ws.Range('D19').Value = '0.0₃0755'
$ws.Range('E19').Value = '  +2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.91'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.03%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  +3.78%  '
$ws.Range('E23').Value = '  +4.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.15'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.29'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.51'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.81'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.36%  '
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  +2.73%  '
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.44'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.72%  '
$ws.Range('E33').Value = '  +3.14%  '
$ws.Range('D34').Value = '1.485.70'
$ws.Range('E34').Value = '  -3.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.66'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.613'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('E37').Value = '  +4.24%  '
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.07'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '71.27'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.07%  '
$ws.Range('E42').Value = '  +6.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.873.87'
$ws.Range('E44').Value = '  +2.99%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.29'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('E46').Value = '  +1.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.73'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +11.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '91.21'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('D49').Value = '0.0₆0111'
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.29'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.105'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.01%  '
